$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking string values are stored as text (matches source cell type)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"

# Apply updated values
$ws.Range('D2').Value = '68.201.55'
$ws.Range('E2').Value = '  -2.24%  '
$ws.Range('D3').Value = '3.593.91'
$ws.Range('E3').Value = '  -2.64%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '626.15'
$ws.Range('E5').Value = '  -6.90%  '
$ws.Range('D6').Value = '155.82'
$ws.Range('E6').Value = '  -3.08%  '
$ws.Range('D7').Value = '3.592.05'
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('E10').Value = '  -3.15%  '
$ws.Range('D11').Value = '6.96'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('D12').Value = '0.434'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').Value = '0.0000225'
$ws.Range('E13').Value = '  -3.57%  '
$ws.Range('D14').Value = '4.212.21'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').Value = '32.09'
$ws.Range('E15').Value = '  -3.46%  '
$ws.Range('D16').Value = '3.621.73'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').Value = '68.290.64'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').Value = '15.68'
$ws.Range('E20').Value = '  -2.75%  '
$ws.Range('D21').Value = '458.44'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').Value = '9.84'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '0.642'
$ws.Range('E23').Value = '  -1.02%  '
$ws.Range('D24').Value = '78.17'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').Value = '3.747.92'
$ws.Range('E25').Value = '  -2.35%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '10.75'
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('E28').Value = '  -8.85%  '
$ws.Range('D29').Value = '8.46'
$ws.Range('E29').Value = '  -7.00%  '
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '26.15'
$ws.Range('E33').Value = '  -2.62%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.92'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('E35').Value = '  -5.04%  '
$ws.Range('D36').Value = '3.604.15'
$ws.Range('E36').Value = '  -2.24%  '
$ws.Range('D37').Value = '6.22'
$ws.Range('E37').Value = '  -4.52%  '
$ws.Range('D38').Value = '8.18'
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('D41').Value = '177.08'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = '5.65'
$ws.Range('E42').Value = '  -7.87%  '
$ws.Range('D43').Value = '2.15'
$ws.Range('E43').Value = '  -5.26%  '
$ws.Range('D44').Value = '0.0880'
$ws.Range('E44').Value = '  -3.26%  '
$ws.Range('D45').Value = '0.902'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '46.13'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '28.71'
$ws.Range('E47').Value = '  +3.67%  '
$ws.Range('E48').Value = '  -5.67%  '
$ws.Range('D49').Value = '7.73'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('E50').Value = '  -7.05%  '
$ws.Range('E51').Value = '  -5.94%  '

# Restore default (unformatted) style so no stray number-format style persists
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D49').Style = "Normal"
